$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "Tesis"
$newRow.Cells.Item(2).Range.Text = "Tesis 101"
$newRow.Cells.Item(3).Range.Text = "En desarrollo"
$newRow.Cells.Item(4).Range.Text = "Profesor Guia"
